# Updates Price (D) and Volume(1h) (E) columns for the cryptos worksheet
# to reflect the latest scraped values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.377.67'
$ws.Range("E2").Value = '  +0.03%  '

# Row 3
$ws.Range("D3").Value = '2.243.64'
$ws.Range("E3").Value = '  -0.10%  '

# Row 4
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.06'
$ws.Range("E5").Value = '  -0.66%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.630'
$ws.Range("E6").Value = '  -0.09%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.79'
$ws.Range("E7").Value = '  -0.18%  '

# Row 8
$ws.Range("E8").Value = '  +0.07%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.622'
$ws.Range("E9").Value = '  -1.77%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.80'
$ws.Range("E10").Value = '  +9.42%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0947'
$ws.Range("E11").Value = '  -0.36%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.27'
$ws.Range("E12").Value = '  +0.64%  '

# Row 13
$ws.Range("E13").Value = '  -1.16%  '

# Row 14
$ws.Range("D14").Value = '2.580.22'
$ws.Range("E14").Value = '  -0.28%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.58'
$ws.Range("E15").Value = '  -1.93%  '

# Row 16
$ws.Range("E16").Value = '  -0.47%  '

# Row 17
$ws.Range("D17").Value = '2.248.43'
$ws.Range("E17").Value = '  +0.57%  '

# Row 18
$ws.Range("D18").Value = '42.172.16'
$ws.Range("E18").Value = '  -0.13%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000102'
$ws.Range("E19").Value = '  +3.82%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.18'
$ws.Range("E20").Value = '  +0.57%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.19'
$ws.Range("E21").Value = '  +1.07%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.23'
$ws.Range("E22").Value = '  +3.35%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.48'
$ws.Range("E23").Value = '  +0.13%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.17'
$ws.Range("E24").Value = '  +32.87%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.41'
$ws.Range("E26").Value = '  +2.36%  '

# Row 27
$ws.Range("E27").Value = '  -2.81%  '

# Row 28
$ws.Range("E28").Value = '  +0.27%  '

# Row 29
$ws.Range("E29").Value = '  +1.48%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.21'
$ws.Range("E30").Value = '  +0.07%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.67'
$ws.Range("E31").Value = '  +0.73%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0829'
$ws.Range("E32").Value = '  -2.78%  '

# Row 33
$ws.Range("E33").Value = '  +0.85%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.71'
$ws.Range("E34").Value = '  -1.13%  '

# Row 35
$ws.Range("E35").Value = '  +0.11%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.27'
$ws.Range("E36").Value = '  +11.40%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.53'
$ws.Range("E37").Value = '  +0.96%  '

# Row 38
$ws.Range("E38").Value = '  +6.72%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.59'
$ws.Range("E39").Value = '  +4.88%  '

# Row 40
$ws.Range("E40").Value = '  -2.03%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.82'
$ws.Range("E41").Value = '  -1.58%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '63.46'
$ws.Range("E42").Value = '  +5.67%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.202'
$ws.Range("E43").Value = '  -0.70%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '108.18'
$ws.Range("E44").Value = '  -7.56%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.78'
$ws.Range("E45").Value = '  +0.41%  '

# Row 46
$ws.Range("E46").Value = '  +1.42%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.996'
$ws.Range("E47").Value = '  -0.06%  '

# Row 48
$ws.Range("E48").Value = '  +2.30%  '

# Row 49
$ws.Range("E49").Value = '  -0.19%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.34'
$ws.Range("E50").Value = '  +5.97%  '

# Row 51
$ws.Range("E51").Value = '  +1.21%  '
